# "cambios de mayo de mayo" - update the SIPOT report period from
# Q4 2021 (Oct-Dec 2021) to Q1 2022 (Jan-Mar 2022) and tweak the header
# row (G3:I3) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")
$ws.Activate()

# --- Header row (row 3): wrap text on the merged G3:I3 cell, and trim
# the border around it down to just a left divider, with a taller row
# to fit the wrapped caption text.
$ws.Rows.Item(3).RowHeight = 66

$g3 = $ws.Range("G3")
$g3.Borders.LineStyle = -4142
$g3.Borders.Item(7).LineStyle = 1
$g3.Borders.Item(7).Weight = 2
$g3.Borders.Item(7).ColorIndex = -4105

$headerG = $ws.Range("G3:I3")
$headerG.WrapText = $true

# --- Row 8 data: shift the reporting period forward one quarter
# (2021 Q4 -> 2022 Q1) and bump the corresponding validation/update dates.
$row8 = $ws.Range("A8")
$row8.Value = 2022

$ws.Range("B8").Value = 44562   # 2022-01-01
$ws.Range("C8").Value = 44651   # 2022-03-31

$ws.Range("AL8").Value = 44659  # 2022-04-08
$ws.Range("AM8").Value = 44659  # 2022-04-08

# --- View state: leave the selection on A8, scrolled back to the top
# of the sheet.
$ws.Range("A8").Select() | Out-Null
